$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '92.241.07'
Set-TextValue 'E2' '  -5.07%  '
Set-TextValue 'D3' '3.318.41'
Set-TextValue 'E3' '  -5.06%  '
Set-TextValue 'E4' '  -0.04%  '
Set-TextValue 'D5' '228.73'
Set-TextValue 'E5' '  -7.41%  '
Set-TextValue 'D6' '614.45'
Set-TextValue 'E6' '  -5.84%  '
Set-TextValue 'D7' '1.35'
Set-TextValue 'E7' '  -5.94%  '
Set-TextValue 'D8' '0.376'
Set-TextValue 'E8' '  -9.53%  '
Set-TextValue 'E9' '  +0.05%  '
Set-TextValue 'D10' '0.909'
Set-TextValue 'E10' '  -10.35%  '
Set-TextValue 'D11' '3.316.17'
Set-TextValue 'E11' '  -5.23%  '
Set-TextValue 'D12' '41.39'
Set-TextValue 'E12' '  -6.24%  '
Set-TextValue 'D13' '0.190'
Set-TextValue 'E13' '  -5.50%  '
Set-TextValue 'D14' '5.92'
Set-TextValue 'E14' '  -4.30%  '
Set-TextValue 'D15' '91.993.63'
Set-TextValue 'E15' '  -5.24%  '
Set-TextValue 'D16' '3.945.22'
Set-TextValue 'E16' '  -4.89%  '
Set-TextValue 'D17' '0.0000239'
Set-TextValue 'E17' '  -6.35%  '
Set-TextValue 'D18' '7.85'
Set-TextValue 'E18' '  -10.08%  '
Set-TextValue 'D19' '3.321.34'
Set-TextValue 'E19' '  -5.04%  '
Set-TextValue 'D20' '17.05'
Set-TextValue 'E20' '  -7.69%  '
Set-TextValue 'D21' '11.06'
Set-TextValue 'E21' '  -8.37%  '
Set-TextValue 'D22' '483.78'
Set-TextValue 'E22' '  -7.19%  '
Set-TextValue 'D23' '3.24'
Set-TextValue 'E23' '  -2.52%  '
Set-TextValue 'D24' '0.435'
Set-TextValue 'E24' '  -13.71%  '
Set-TextValue 'E25' '  -9.96%  '
Set-TextValue 'D26' '6.06'
Set-TextValue 'E26' '  -10.39%  '
Set-TextValue 'D27' '88.88'
Set-TextValue 'E27' '  -8.08%  '
Set-TextValue 'D28' '3.499.40'
Set-TextValue 'E28' '  -4.96%  '
Set-TextValue 'D29' '11.45'
Set-TextValue 'E29' '  -9.59%  '
Set-TextValue 'E30' '  +0.07%  '
Set-TextValue 'D31' '10.99'
Set-TextValue 'E31' '  -9.25%  '
Set-TextValue 'D32' '0.134'
Set-TextValue 'E32' '  -4.59%  '
Set-TextValue 'D33' '2.62'
Set-TextValue 'E33' '  -5.81%  '
Set-TextValue 'D34' '0.998'
Set-TextValue 'E34' '  +0.08%  '
Set-TextValue 'D35' '0.170'
Set-TextValue 'E35' '  -9.15%  '
Set-TextValue 'D36' '27.91'
Set-TextValue 'E36' '  -11.00%  '
Set-TextValue 'D37' '0.520'
Set-TextValue 'E37' '  -11.27%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D38' '518.87'
Set-TextValue 'E38' '  +0.99%  '
$ws.Range('B39').Value = 'USDe'
$ws.Range('C39').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D39' '1.00'
Set-TextValue 'E39' '  -0.05%  '
Set-TextValue 'D40' '7.23'
Set-TextValue 'E40' '  -8.54%  '
Set-TextValue 'D41' '0.145'
Set-TextValue 'E41' '  -6.23%  '
Set-TextValue 'D42' '1.33'
Set-TextValue 'E42' '  -9.98%  '
Set-TextValue 'D43' '0.870'
Set-TextValue 'E43' '  -4.81%  '
Set-TextValue 'D44' '23.99'
Set-TextValue 'E44' '  -1.27%  '
$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D45' '1.64'
Set-TextValue 'E45' '  -4.27%  '
$ws.Range('B46').Value = 'MantraDAO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextValue 'D46' '3.55'
Set-TextValue 'E46' '  -0.98%  '
Set-TextValue 'D47' '0.0396'
Set-TextValue 'E47' '  -6.74%  '
Set-TextValue 'D48' '5.28'
Set-TextValue 'E48' '  -6.24%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D49' '2.08'
Set-TextValue 'E49' '  -6.10%  '
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D50' '51.80'
Set-TextValue 'E50' '  -4.40%  '
Set-TextValue 'E51' '  -7.72%  '
